# Append a page break followed by a new "Appendix" section:
#   - "Додаток А" (centered title)
#   - blank centered paragraph
#   - "Посилання на сайт: https://tervild.github.io/webpractical/" (left aligned)
#
# wdAlignParagraphLeft   = 0
# wdAlignParagraphCenter = 1
# wdPageBreak            = 7

$d = $word.ActiveDocument

$lastPara = $d.Paragraphs.Last
$tail = $lastPara.Range
$tail.Collapse(0)

# --- New paragraph that holds only the page break -----------------------
$tail.InsertParagraphAfter()
$breakPara = $d.Paragraphs.Last
$breakRange = $breakPara.Range
$breakRange.Collapse(0)
$breakRange.Text = [char]12

# --- "Додаток А" title, centered -----------------------------------------
$breakRange2 = $breakPara.Range
$breakRange2.Collapse(0)
$breakRange2.InsertParagraphAfter()
$titlePara = $d.Paragraphs.Last
$titlePara.Alignment = 1
$titleRange = $titlePara.Range
$titleRange.Collapse(0)
$titleRange.Font.Name = "Times New Roman"
$titleRange.Font.NameFarEast = "Times New Roman"
$titleRange.Text = "Додаток"
$titleRange2 = $titlePara.Range
$titleRange2.Collapse(0)
$titleRange2.Font.Name = "Times New Roman"
$titleRange2.Font.NameFarEast = "Times New Roman"
$titleRange2.Text = " А"

# --- Blank centered paragraph ---------------------------------------------
$titleRange3 = $titlePara.Range
$titleRange3.Collapse(0)
$titleRange3.InsertParagraphAfter()
$blankPara = $d.Paragraphs.Last
$blankPara.Alignment = 1

# --- "Посилання на сайт: ..." paragraph, left aligned ---------------------
$blankRange = $blankPara.Range
$blankRange.Collapse(0)
$blankRange.InsertParagraphAfter()
$linkPara = $d.Paragraphs.Last
$linkPara.Alignment = 0
$linkRange = $linkPara.Range
$linkRange.Collapse(0)
$linkRange.Font.Name = "Times New Roman"
$linkRange.Font.NameFarEast = "Times New Roman"
$linkRange.Text = "Посилання"

$linkRange2 = $linkPara.Range
$linkRange2.Collapse(0)
$linkRange2.Font.Name = "Times New Roman"
$linkRange2.Font.NameFarEast = "Times New Roman"
$linkRange2.Text = " "

$linkRange3 = $linkPara.Range
$linkRange3.Collapse(0)
$linkRange3.Font.Name = "Times New Roman"
$linkRange3.Font.NameFarEast = "Times New Roman"
$linkRange3.Text = "на"

$linkRange4 = $linkPara.Range
$linkRange4.Collapse(0)
$linkRange4.Font.Name = "Times New Roman"
$linkRange4.Font.NameFarEast = "Times New Roman"
$linkRange4.Text = " "

$linkRange5 = $linkPara.Range
$linkRange5.Collapse(0)
$linkRange5.Font.Name = "Times New Roman"
$linkRange5.Font.NameFarEast = "Times New Roman"
$linkRange5.Text = "сайт"

$linkRange6 = $linkPara.Range
$linkRange6.Collapse(0)
$linkRange6.Font.Name = "Times New Roman"
$linkRange6.Font.NameFarEast = "Times New Roman"
$linkRange6.Text = ": https://tervild.github.io/webpractical/"

Write-Output ("Paragraphs now: " + $d.Paragraphs.Count)
